# "fixed template, fixed legal condition"
#
# 1) TTC formula (row 30) used a hard-coded 1.2 multiplier; make it follow the
#    VAT rate cell ($C$8, expressed as a percentage) instead.
# 2) The legal-conditions block (DPA annex + CGS hyperlinks/dates, rows
#    43-47) referenced an outdated contract revision, so it is removed
#    entirely; everything below it shifts up by 4 rows.
# 3) The print area is shrunk to match the new, shorter sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quote")

# --- 1) Generalise the TTC-budget formula to use the VAT-rate cell ---------
$ws.Range("J30").Formula = '=J29*(1+$C$8/100)'

# --- 2) Drop the obsolete DPA/CGS legal-reference rows ----------------------
# A43 used to hold "L'Annexe <<Traitement de Donnees ...>>" text; it becomes
# a bare styled cell once the block below it (old rows 44-47, the hyperlink /
# date lines) is deleted.
$ws.Range("A43").ClearContents()
$ws.Rows("44:47").Delete()

# --- 3) Shrink the print area to the new sheet extent -----------------------
$ws.PageSetup.PrintArea = '$A$1:$J$81'

# --- cosmetic: restore selection/scroll near where the edit happened -------
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("A44").Select()
